# Loginshizzle toegevoegd, woei! Uurtjes bijgewerkt.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New hour entries for Mr. Pink (E/F columns) and Mr. White (B/C columns).
# Dates are stored as Excel date serials; styles (s="1" date, s="4" number)
# already exist on these cells so just setting .Value keeps formatting.
$ws.Range("E6").Value = 41062
$ws.Range("F6").Value = 6.5

$ws.Range("E7").Value = 41063
$ws.Range("F7").Value = 2

$ws.Range("B8").Value = 41062
$ws.Range("C8").Value = 6.5

# Recalculate the SUM()/derived formulas with the new hours.
$excel.Calculate()

# Update the active selection to match the recorded view state.
$ws.Range("R12").Select()
